$d = $word.ActiveDocument

# --- Edit 1: split the run "left-bottom" into two runs: "left-" and "middle" ---
# Locate the original "left-bottom" text.
$rng = $d.Content
$found = $rng.Find.Execute("left-bottom", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # "left-bottom" -> "left-" (5 chars) + "bottom" (6 chars). Grab the trailing
    # "bottom" substring so it can be turned into its own run carrying "middle".
    $tailStart = $rng.End - 6
    $tailEnd = $rng.End
    $tailRng = $d.Range($tailStart, $tailEnd)

    # Touching the font (clearing to "automatic") forces Word to split this
    # substring into its own run (distinct rPr) instead of merging it back
    # into the "left-" run, and leaves the original "left-" run/run-properties
    # untouched.
    $tailRng.Font.Color = -16777216

    # Replace the "bottom" text with "middle" within that now-separate run.
    $tailRng.Text = "middle"
}

# --- Edit 2: give the section a footer distance (w:footer="720") ---
$d.PageSetup.FooterDistance = 36
